$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 4: correct the mistyped "Fecha Pago" date (.12.15.2022 -> .15.12.2022) ---
$ws.Range("L4").Value = ".15.12.2022"

# --- Fill column A (Fecha) for the new rows 5, 6, 7 ---
$ws.Range("A5").Font.Underline = $false
$ws.Range("A5").Value = ".06.02.2021"
$ws.Range("A6").Value = ".06.02.2021"
$ws.Range("A7").Value = ".13.07.2021"

# --- Fill column L (Fecha Pago) for the new rows 5, 6, 7 ---
$ws.Range("L5").Value = ".11.02.2021"
$ws.Range("L6").Value = ".11.02.2021"
$ws.Range("L7").Value = ".15.07.2021"

# --- Row 5: rest of the columns (mirrors row 2's record) ---
$ws.Range("B5").Value = "E220106-01-01  "
$ws.Range("C5").Value = 5409
$ws.Range("D5").Value = "PIURA               "
$ws.Range("E5").Value = "AAAAA"
$ws.Range("F5").Value = "AAAA"
$ws.Range("G5").Value = "CEDULA              "
$ws.Range("H5").Value = 1243333333
$ws.Range("I5").Value = 500.5
$ws.Range("I5").NumberFormat = "#,##0.00"
$ws.Range("J5").Value = "PAGADO    "
$ws.Range("K5").Value = "PIURA                    "
$ws.Range("M5").Value = 0.50050925925925926
$ws.Range("M5").NumberFormat = "h:mm:ss"
$ws.Range("N5").Value = 99999999
$ws.Range("O5").Value = 797453
$ws.Range("P5").Value = "                    "
$ws.Range("Q5").Value = "                    "
$ws.Range("R5").Value = "Estadia 1           "
$ws.Range("S5").Value = "BHA            "

# --- Row 6: rest of the columns (mirrors row 3's record) ---
$ws.Range("B6").Value = "E220106-01-01  "
$ws.Range("C6").Value = 5410
$ws.Range("D6").Value = "PIURA               "
$ws.Range("E6").Value = "BBBBB"
$ws.Range("F6").Value = "BBBB"
$ws.Range("G6").Value = "CEDULA              "
$ws.Range("H6").Value = 9876543
$ws.Range("I6").Value = 5.5
$ws.Range("I6").NumberFormat = "#,##0.00"
$ws.Range("J6").Value = "DEVUELTO  "
$ws.Range("K6").Value = "PIURA                    "
$ws.Range("M6").Value = "        "
$ws.Range("N6").Value = 918881831
$ws.Range("O6").Value = 797453
$ws.Range("P6").Value = "                    "
$ws.Range("Q6").Value = "                    "
$ws.Range("R6").Value = "Estadia 1           "
$ws.Range("S6").Value = "BHA            "

# --- Row 7: rest of the columns (mirrors row 4's record) ---
$ws.Range("B7").Value = "E220106-01-01  "
$ws.Range("C7").Value = 5410
$ws.Range("D7").Value = "PIURA               "
$ws.Range("E7").Value = "BBBBB"
$ws.Range("F7").Value = "BBBB"
$ws.Range("G7").Value = "CEDULA              "
$ws.Range("H7").Value = 9876543
$ws.Range("I7").Value = 1500.9
$ws.Range("I7").NumberFormat = "#,##0.00"
$ws.Range("J7").Value = "DEVUELTO  "
$ws.Range("K7").Value = "PIURA                    "
$ws.Range("M7").Value = "        "
$ws.Range("N7").Value = 918881831
$ws.Range("O7").Value = 797453
$ws.Range("P7").Value = "                    "
$ws.Range("Q7").Value = "                    "
$ws.Range("R7").Value = "Estadia 1           "
$ws.Range("S7").Value = "BHA            "

# --- Final selection, as recorded in the saved view state ---
$ws.Range("L8").Select()
